# "Generate Report for handback"
#
# The handback-status report is refreshed now that a handback has
# occurred: the status moves from "Not localized" to "Handed back", the
# "Latest ..." column headers are renamed to "Correspond ..." (and the
# "Latest Target File" header simplifies to "Target File"), and the
# per-language sheets gain a populated "Correspond Handoff File" column
# (C) carrying the handoff .xlf file name/hyperlink that was previously
# left blank.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Status: "Not localized" -> "Handed back" on every sheet.
# ---------------------------------------------------------------
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B2").Value = "Handed back"
    if ($sheetName -eq "Overview") {
        $ws.Range("C2").Value = "Handed back"
    }
}

# ---------------------------------------------------------------
# 2. Per-language sheets: rename headers + fill in the handoff file.
# ---------------------------------------------------------------
$langs = @{
    "zh-cn" = @{
        HandoffFile   = "fb725adb-35ea-46a1-9779-7d7ffb4de3e1.e7acdc3f3dd1279889781c8d2132960756381a10.zh-cn.xlf"
        HandoffTarget = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0b44c2c5a6a59fc2c07ae027d3c4a666b90cb283/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/fb725adb-35ea-46a1-9779-7d7ffb4de3e1.e7acdc3f3dd1279889781c8d2132960756381a10.zh-cn.xlf"
    }
    "de-de" = @{
        HandoffFile   = "fb725adb-35ea-46a1-9779-7d7ffb4de3e1.e7acdc3f3dd1279889781c8d2132960756381a10.de-de.xlf"
        HandoffTarget = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/847ace712fff9f4fc90ab4b09e802a09f556a0d0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/fb725adb-35ea-46a1-9779-7d7ffb4de3e1.e7acdc3f3dd1279889781c8d2132960756381a10.de-de.xlf"
    }
}

foreach ($lang in $langs.Keys) {
    $ws = $wb.Worksheets.Item($lang)
    $info = $langs[$lang]

    # Rename the "Latest ..." headers to "Correspond ..." (and simplify
    # "Latest Target File" to just "Target File").
    $ws.Range("C1").Value = "Correspond Handoff File"
    $ws.Range("D1").Value = "Correspond Handoff Datetime"
    $ws.Range("E1").Value = "Target File"
    $ws.Range("F1").Value = "Correspond Handback File"
    $ws.Range("G1").Value = "Correspond Handback DateTime"

    # Fill in the previously-empty "Correspond Handoff File" cell with
    # the handoff .xlf file name, styled/linked like the other file-name
    # hyperlink cells (A2, E2, F2).
    $ws.Range("C2").Value = $info.HandoffFile
    $ws.Hyperlinks.Add($ws.Range("C2"), $info.HandoffTarget, [Type]::Missing, [Type]::Missing, $info.HandoffFile) | Out-Null
    $ws.Range("C2").Font.Underline = $true
    $ws.Range("C2").Font.Color = $ws.Range("F2").Font.Color
}
